# Applies:
#  1. Refresh the cached "datetimeFigureOut" footer-date field text (12/2/22 -> 2/19/24)
#     on the slide master, every slide layout, and the notes master.
#  2. Tweak the wording on the "FollowUp intro" slide's instructions textbox and
#     let the shape's auto-fit height follow the now-shorter text.

$p = $ppt.ActivePresentation

function Update-DatePlaceholder {
    param($shapes)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.PlaceholderFormat.Type -eq 16) {
            # ppPlaceholderDate
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "12/2/22") {
                $tr.Characters(1, $tr.Text.Length).Text = "2/19/24"
            }
        }
    }
}

# -- Slide master date placeholder --
Update-DatePlaceholder $p.SlideMaster.Shapes

# -- Every slide layout's date placeholder --
$layouts = $p.SlideMaster.CustomLayouts
for ($L = 1; $L -le $layouts.Count; $L++) {
    Update-DatePlaceholder $layouts.Item($L).Shapes
}

# -- Notes master date placeholder --
Update-DatePlaceholder $p.NotesMaster.Shapes

# -- Slide 1: "FollowUp intro" instructions textbox --
$slide = $p.Slides.Item(1)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        $tr = $shp.TextFrame.TextRange
        $full = $tr.Text

        $oldIntro = "Congratulations! You have successfully completed all trials. Thank you for your participation in our study. "
        $newIntro = "Congratulations! You have successfully completed all trials. "
        $oldAsk = "Please answer these short follow up questions. Press the F key to begin."
        $newAsk = "Please answer the following questions. Click the mouse to begin."

        $idx1 = $full.IndexOf($oldIntro)
        if ($idx1 -ge 0) {
            $tr.Characters($idx1 + 1, $oldIntro.Length).Text = $newIntro
        }

        $full = $tr.Text
        $idx2 = $full.IndexOf($oldAsk)
        if ($idx2 -ge 0) {
            $tr.Characters($idx2 + 1, $oldAsk.Length).Text = $newAsk
        }
    }
}
